# DTT Assessment Hour Log - "Did a TestRun not done yet"
# The test run on 13-12-2023 (row 8) actually only took 10 hours instead of 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the logged amount of hours for that day (B8: 12 -> 10).
# B30 (Total amount of hours) holds =SUMIF(E4:E28,"<>x",B4:B28) and will
# recalculate automatically from this change.
$ws.Range("B8").Value = 10

# Move the current selection to B25 (no frozen/scrolled top-left cell anymore).
$ws.Range("B25").Select()
